$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rights_and_functions")

# Insert a new row above row 36 (shifts existing rows 36-46 down to 37-47)
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with SELECT / cds2db_user (same right pattern as row 32)
$ws.Range("I36").Value = "SELECT"
$ws.Range("J36").Value = "cds2db_user"

# Change right for frontend_user: the db2frontend_in row (now row 43 after the insert)
# previously listed db2dataprocessor_user in column J - fix it to db2frontend_user
$ws.Range("J43").Value = "db2frontend_user"

# The row insert shifts cell data but this engine does not auto-relocate cell comments,
# so move the three comments that sat in the shifted rows to their new cell addresses.
$txt = "Autor:`nWie 30 und 31"

$ws.Range("K37").Comment.Delete()
$ws.Range("K38").AddComment($txt)

$ws.Range("K39").Comment.Delete()
$ws.Range("K40").AddComment($txt)

$ws.Range("K42").Comment.Delete()
$ws.Range("K43").AddComment($txt)

# Update the selected cell to match the target workbook view state
$ws.Range("J36").Select()

# Configure page setup for printing (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
